# Applies the weekly data refresh: for each data row (2..50) the values in
# columns D (Fecha) and L:T (Calidad .. Kg/unidad) are re-shuffled among the
# rows according to the mapping below (columns A:C and E:K are untouched).
# new_row[r], cols D,L:T  =  old_row[ $rowMap[r] ], cols D,L:T

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 50

$rowMap = @{
    2 = 2
    3 = 3
    4 = 17
    5 = 32
    6 = 33
    7 = 5
    8 = 6
    9 = 7
    10 = 14
    11 = 23
    12 = 24
    13 = 11
    14 = 34
    15 = 35
    16 = 37
    17 = 36
    18 = 47
    19 = 48
    20 = 21
    21 = 10
    22 = 42
    23 = 41
    24 = 30
    25 = 12
    26 = 43
    27 = 44
    28 = 18
    29 = 16
    30 = 45
    31 = 9
    32 = 31
    33 = 4
    34 = 28
    35 = 29
    36 = 25
    37 = 49
    38 = 50
    39 = 26
    40 = 27
    41 = 15
    42 = 40
    43 = 39
    44 = 8
    45 = 13
    46 = 46
    47 = 38
    48 = 19
    49 = 20
    50 = 22
}

# snapshot the "before" values for the columns that move (D and L..T)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# write the new values back according to the row map
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $snapshot[$rowMap[$r]]
    $ws.Cells.Item($r, 4).Value = $src.D
    $ws.Cells.Item($r, 12).Value = $src.L
    $ws.Cells.Item($r, 13).Value = $src.M
    $ws.Cells.Item($r, 14).Value = $src.N
    $ws.Cells.Item($r, 15).Value = $src.O
    $ws.Cells.Item($r, 16).Value = $src.P
    $ws.Cells.Item($r, 17).Value = $src.Q
    $ws.Cells.Item($r, 18).Value = $src.R
    $ws.Cells.Item($r, 19).Value = $src.S
    $ws.Cells.Item($r, 20).Value = $src.T
}
